# Insert a new weekly price record as row 20, pushing all existing
# records (old rows 20-130) down by one (to 21-131).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(20).Insert()

# The newly inserted row 20 is blank; seed it with the same row "shape"
# as the record that used to occupy row 20 (now shifted to row 21), then
# overwrite just the fields that differ for this new entry.
$ws.Range("A21:R21").Copy()
$ws.Range("A20:R20").PasteSpecial()

$ws.Range("D20").Value = 45250
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 30000
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 1200
